# Powerpoint writer: consolidate text run nodes.
#
# The original XML had each word and the space that follows it living in
# their own separate <a:r> run (e.g. "Here" / " " / "is" / " " / "a" / ...).
# The target XML merges each word with its trailing space into a single
# run (e.g. "Here " / "is " / "a " / ...), which is what PowerPoint itself
# does when you simply retype/replace text through the UI - it keeps the
# existing runs it can match a common prefix against and only rewrites the
# minimal trailing piece, naturally consolidating the "word + following
# space" pairs into one run each.
#
# We reproduce that consolidation here by re-assigning the text of
# "Characters" sub-ranges that exactly span [word][following space] for
# each word/space pair except the very last word, which is left untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# Slide title: "Here is a single header" -> runs "Here ","is ","a ","single ","header"
# ---------------------------------------------------------------------
$titleTr = $s.Shapes.Item(1).TextFrame.TextRange

$titleTr.Characters(1, 5).Text  = "Here "
$titleTr.Characters(6, 3).Text  = "is "
$titleTr.Characters(9, 2).Text  = "a "
$titleTr.Characters(11, 7).Text = "single "
# "header" (chars 18-23) is already its own run and stays as-is.

# ---------------------------------------------------------------------
# Speaker notes: "and here are some notes" -> runs "and ","here ","are ","some ","notes"
# ---------------------------------------------------------------------
$notesTr = $s.NotesPage.Shapes.Item(2).TextFrame.TextRange
$notesTr.Text = "and here are some notes"
